{"js": "// Update benchmark stats table: each table row has exactly one cell.\n// Map of (0-based) row index -> new cell text.\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"304\",\n  4: \"0.00003\",\n  6: \"0.00016\",\n  7: \"0.00006\",\n  8: \"0.00020\",\n  9: \"0.00023\",\n  10: \"0.00025\",\n  11: \"0.04935\",\n  43: \"99.95\",\n  44: \"0.05\",\n  45: \"99\",\n};\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst maxIndex = Math.max(...Object.keys(updates).map(Number));\nfor (let i = 0; i <= maxIndex && i < rows.items.length; i++) {\n  if (!(i in updates)) continue;\n  const row = rows.items[i];\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n  // Each row in this table has a single cell holding the value.\n  const cell = cells.items[0];\n  cell.value = updates[i];\n}\n\nawait context.sync();\n", "ps1": "# Update benchmark stats table: each table row has exactly one cell\n# (column index 1). Word COM table cell/row indices are 1-based.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"0M\"\n$t.Cell(2, 1).Range.Text = \"0M\"\n$t.Cell(3, 1).Range.Text = \"0M\"\n$t.Cell(4, 1).Range.Text = \"304\"\n$t.Cell(5, 1).Range.Text = \"0.00003\"\n$t.Cell(7, 1).Range.Text = \"0.00016\"\n$t.Cell(8, 1).Range.Text = \"0.00006\"\n$t.Cell(9, 1).Range.Text = \"0.00020\"\n$t.Cell(10, 1).Range.Text = \"0.00023\"\n$t.Cell(11, 1).Range.Text = \"0.00025\"\n$t.Cell(12, 1).Range.Text = \"0.04935\"\n$t.Cell(44, 1).Range.Text = \"99.95\"\n$t.Cell(45, 1).Range.Text = \"0.05\"\n$t.Cell(46, 1).Range.Text = \"99\"\n"}
